$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E7").Value = 11.6844
$ws.Range("C10").Value = -12.9504
$ws.Range("C12").Value = -14.29500000000001
$ws.Range("D13").Value = -7.988300000000001
$ws.Range("C18").Value = -14.1037
$ws.Range("E20").Value = 13.01249999999998
